$wb = $excel.ActiveWorkbook

# 1) Rename "Requested quantity" header to the new metric-specific names.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add a new "PO Forecast" sheet at the end of the workbook, copying the
#    existing "Monthly Trend" sheet so the new sheet inherits the same
#    look (bordered/bold header row, date-formatted first column) before
#    we overwrite its content with the forecast data.
$wsMonthly.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"

# Header row. A1/B1 already carry the bold/bordered header style inherited
# from the "Monthly Trend" copy; C1/D1 are new cells, so copy that same
# header format onto them (reuses the existing style instead of minting a
# new one).
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$wsForecast.Range("A1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @()
$data += ,@(45039.99999999999, 0, -93.70740575670885, 47.29845428249723)
$data += ,@(45123.99999999999, 62, -15.56728951905131, 138.8711215173358)
$data += ,@(45130.99999999999, 69, -10.08916729526367, 147.2790082877599)
$data += ,@(45144.99999999999, 83, 8.903139697279558, 164.9742155026427)
$data += ,@(45151.99999999999, 90, 18.86778423287475, 159.9234153651528)
$data += ,@(45158.99999999999, 98, 29.5035289845801, 173.7733561577026)
$data += ,@(45193.99999999999, 133, 56.97982188860703, 206.4444332479123)
$data += ,@(45221.99999999999, 161, 84.52795427722299, 231.6770424609488)
$data += ,@(45235.99999999999, 175, 96.79763872079714, 246.7712213437229)
$data += ,@(45242.99999999999, 182, 104.7568477764057, 249.514548028052)
$data += ,@(45249.99999999999, 189, 116.5422265701818, 263.978092600337)
$data += ,@(45277.99999999999, 218, 142.7786120145434, 284.3261177671708)
$data += ,@(45284.99999999999, 225, 150.9020580365568, 295.0068830390557)
$data += ,@(45291.99999999999, 232, 159.0681513728983, 307.4885023466886)
$data += ,@(45298.99999999999, 239, 168.5116868468086, 312.7580797579206)
$data += ,@(45305.99999999999, 246, 172.1951055461445, 321.6157602819924)
$data += ,@(45312.99999999999, 253, 180.69997651063, 327.9467478553659)
$data += ,@(45319.99999999999, 260, 194.9586571227886, 339.8468653247016)
$data += ,@(45326.99999999999, 267, 185.5612609500765, 341.3423242607347)
$data += ,@(45333.99999999999, 274, 200.8405042246795, 345.7430255807378)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$wsForecast.Range("A1").Select() | Out-Null
